$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (A:T), reflecting the expanded
# sending-cluster x target-cluster grid (ECs added as a sending cluster,
# per Dr Hou's advice) for the Slit2 -> App ligand-receptor pair.
$data = @(
    @("ECs", "Slit2", "App", "ECs", 2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 82.83048866666665, 248.491466, 0.3167437020391103, 0.3167437020391103, 2.769299337755555, 24.9236940398, 0.005553339262809897, 0.005553339262809897),
    @("ECs", "Slit2", "App", "FAPs", 2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 152.851481, 458.554443, 0.5845039034954311, 0.5845039034954312, 5.110334514766667, 45.9930106329, 0.01024787061479135, 0.01024787061479135),
    @("ECs", "Slit2", "App", "sCs", 2, 0.6666666666666666, 0.03343333333333334, 0.1003, 0.01753259568243662, 0.01753259568243662, 3, 1, 25.824378, 77.47313399999999, 0.09875239446545848, 0.0987523944654585, 0.8633950377999999, 7.770555340199999, 0.001731385804835376, 0.001731385804835376),
    @("FAPs", "Slit2", "App", "ECs", 3, 1, 0.37963, 1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 82.83048866666665, 248.491466, 0.3167437020391103, 0.3167437020391103, 31.44493841252666, 283.0044457127399, 0.0630572537689089, 0.06305725376890892),
    @("FAPs", "Slit2", "App", "FAPs", 3, 1, 0.37963, 1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 152.851481, 458.554443, 0.5845039034954311, 0.5845039034954312, 58.02700773203, 522.2430695882699, 0.1163628849898277, 0.1163628849898277),
    @("FAPs", "Slit2", "App", "sCs", 3, 1, 0.37963, 1.13889, 0.1990797397484571, 0.1990797397484571, 3, 1, 25.824378, 77.47313399999999, 0.09875239446545848, 0.0987523944654585, 9.803708620139998, 88.23337758125999, 0.01965960098972045, 0.01965960098972045),
    @("sCs", "Slit2", "App", "ECs", 3, 1, 1.493861, 4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 82.83048866666665, 248.491466, 0.3167437020391103, 0.3167437020391103, 123.7372366300753, 1113.635129670678, 0.2481331090073915, 0.2481331090073916),
    @("sCs", "Slit2", "App", "FAPs", 3, 1, 1.493861, 4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 152.851481, 458.554443, 0.5845039034954311, 0.5845039034954312, 228.338866258141, 2055.049796323269, 0.4578931478908121, 0.4578931478908122),
    @("sCs", "Slit2", "App", "sCs", 3, 1, 1.493861, 4.481583000000001, 0.7833876645691064, 0.7833876645691064, 3, 1, 25.824378, 77.47313399999999, 0.09875239446545848, 0.0987523944654585, 38.578031143458, 347.202280291122, 0.07736140767090266, 0.07736140767090267)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
